$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.642.78'
$ws.Range('E2').Value = '  +13.30%  '
$ws.Range('D3').Value = '1.829.24'
$ws.Range('E3').Value = '  +9.34%  '
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '230.96'
$ws.Range('E5').Value = '  +5.13%  '
$ws.Range('D6').Value = '0.550'
$ws.Range('E6').Value = '  +4.12%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '31.79'
$ws.Range('E8').Value = '  +7.41%  '
$ws.Range('D9').Value = '46.99'
$ws.Range('E9').Value = '  +6.53%  '
$ws.Range('E10').Value = '  +7.28%  '
$ws.Range('D11').Value = '0.0673'
$ws.Range('E11').Value = '  +4.83%  '
$ws.Range('D12').Value = '0.0931'
$ws.Range('E12').Value = '  +2.85%  '
$ws.Range('D13').Value = '2.092.03'
$ws.Range('E13').Value = '  +9.36%  '
$ws.Range('D14').Value = '1.825.88'
$ws.Range('E14').Value = '  +9.32%  '
$ws.Range('D15').Value = '0.651'
$ws.Range('E15').Value = '  +6.12%  '
$ws.Range('D16').Value = '34.576.05'
$ws.Range('E16').Value = '  +13.08%  '
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '4.28'
$ws.Range('E18').Value = '  +6.78%  '
$ws.Range('E19').Value = '  +5.08%  '
$ws.Range('D20').Value = '259.80'
$ws.Range('E20').Value = '  +6.87%  '
$ws.Range('E21').Value = '  +4.38%  '
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').Value = '10.57'
$ws.Range('E23').Value = '  +5.66%  '
$ws.Range('D24').Value = '4.36'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').Value = '2.21'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').Value = '158.37'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '16.72'
$ws.Range('E27').Value = '  +5.37%  '
$ws.Range('D28').Value = '7.13'
$ws.Range('E28').Value = '  +6.73%  '
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').Value = '0.995'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('E31').Value = '  +12.30%  '
$ws.Range('D32').Value = '0.0520'
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('E33').Value = '  +6.50%  '
$ws.Range('D34').Value = '3.58'
$ws.Range('E34').Value = '  +8.81%  '
$ws.Range('D35').Value = '1.554.08'
$ws.Range('E35').Value = '  +4.43%  '
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('E37').Value = '  +6.06%  '
$ws.Range('B38').Value = 'MinaProtocolToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D38').Value = '1.29'
$ws.Range('E38').Value = '  +212.83%  '
$ws.Range('E39').Value = '  +6.85%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '0.634'
$ws.Range('E40').Value = '  +5.82%  '
$ws.Range('D41').Value = '85.05'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').Value = '2.82'
$ws.Range('E42').Value = '  +5.56%  '
$ws.Range('D43').Value = '0.921'
$ws.Range('E43').Value = '  +9.61%  '
$ws.Range('D44').Value = '2.33'
$ws.Range('E44').Value = '  +1.62%  '
$ws.Range('E45').Value = '  +10.22%  '
$ws.Range('D46').Value = '0.0525'
$ws.Range('E46').Value = '  +5.38%  '
$ws.Range('E47').Value = '  +5.48%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.991.85'
$ws.Range('E48').Value = '  +10.19%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '12.42'
$ws.Range('E49').Value = '  +27.07%  '
$ws.Range('E50').Value = '  +5.30%  '
$ws.Range('D51').Value = '53.22'
$ws.Range('E51').Value = '  +3.82%  '
